$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the existing "sum" header cell (G1) onto the
# new "Save" header cell (H1), so it matches the bold/centered/bordered look
# of the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header label and the new data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
